$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (ResponseId 7): three answers were re-scored to the "Bad 1" choice ---
$ws.Range("F8").Value = [char]0x00A0 + "Bad 1"
$ws.Range("H8").Value = [char]0x00A0 + "Bad 1"
$ws.Range("K8").Value = [char]0x00A0 + "Bad 1"

# --- Row 10 (ResponseId 9): rater identity corrected back to arda.aydin@uzh.ch ---
$ws.Range("C10").Value = "arda.aydin@uzh.ch"
$ws.Range("D10").Value = "arda.aydin@uzh.ch"

# that rater's rows use the default style/row-height (no special 12pt font, no 15.75 row height)
$ws.Range("C10:D10").Style = "Normal"
$ws.Rows.Item(10).AutoFit()

# --- Update the active selection to reflect where the editor ended up ---
[void]$ws.Range("D10").Select()
